# 2DES BCD AULAS REORGANIZADAS
# Insert two new "aula" (class) columns (AU and AV) for the BCD subject
# block, filling them with "P" (presente) for every student row, mirroring
# the existing AT/AW columns. Row 13 is a hidden/withdrawn student row and
# is intentionally left untouched (it has no marks in this column range).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 3; $row -le 30; $row++) {
    if ($row -eq 13) {
        continue
    }
    $ws.Range("AU$row").Value = "P"
    $ws.Range("AV$row").Value = "P"
}

# Restore the (new) active selection recorded for this edit.
$ws.Range("AW15").Select()
